$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rebuild the data rows (1-14) from scratch, in document order, so the
# shared-strings table is re-created the same way Excel lays it out when it
# walks the sheet top-to-bottom / left-to-right on save. A new
# "N3XeS.CSharp.Common.Exceptions" class (ReadOnlyException with its four
# constructors) was inserted into the report, and every other data
# row/metric was refreshed from the regenerated Code Metrics run.
# ---------------------------------------------------------------------------
$ws.Range("A1:J14").ClearContents()

# Row 1 - headers (unchanged, re-entered so they keep their original order)
$ws.Range("A1").Value = "Scope"
$ws.Range("B1").Value = "Project"
$ws.Range("C1").Value = "Namespace"
$ws.Range("D1").Value = "Type"
$ws.Range("E1").Value = "Member"
$ws.Range("F1").Value = "Maintainability Index"
$ws.Range("G1").Value = "Cyclomatic Complexity"
$ws.Range("H1").Value = "Depth of Inheritance"
$ws.Range("I1").Value = "Class Coupling"
$ws.Range("J1").Value = "Lines of Code"

# Row 2 - Project summary totals
$ws.Range("A2").Value = "Project"
$ws.Range("B2").Value = "Source\N3XeS.CSharp.Common (Release)"
$ws.Range("C2").Value = " "
$ws.Range("D2").Value = " "
$ws.Range("E2").Value = " "
$ws.Range("F2").Value = 97
$ws.Range("G2").Value = 145
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 21
$ws.Range("J2").Value = 191

# Row 3 - N3XeS.CSharp.Common.Exceptions (namespace)
$ws.Range("A3").Value = "Namespace"
$ws.Range("B3").Value = "Source\N3XeS.CSharp.Common (Release)"
$ws.Range("C3").Value = "N3XeS.CSharp.Common.Exceptions"
$ws.Range("D3").Value = " "
$ws.Range("E3").Value = " "
$ws.Range("F3").Value = 93
$ws.Range("G3").Value = 4
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 4

# Row 4 - ReadOnlyException (type)
$ws.Range("A4").Value = "Type"
$ws.Range("B4").Value = "Source\N3XeS.CSharp.Common (Release)"
$ws.Range("C4").Value = "N3XeS.CSharp.Common.Exceptions"
$ws.Range("D4").Value = "ReadOnlyException"
$ws.Range("E4").Value = " "
$ws.Range("F4").Value = 93
$ws.Range("G4").Value = 4
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 4

# Row 5 - ReadOnlyException()
$ws.Range("A5").Value = "Member"
$ws.Range("B5").Value = "Source\N3XeS.CSharp.Common (Release)"
$ws.Range("C5").Value = "N3XeS.CSharp.Common.Exceptions"
$ws.Range("D5").Value = "ReadOnlyException"
$ws.Range("E5").Value = "ReadOnlyException()"
$ws.Range("F5").Value = 98
$ws.Range("G5").Value = 1
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 1

# Row 6 - ReadOnlyException(SerializationInfo, StreamingContext)
$ws.Range("A6").Value = "Member"
$ws.Range("B6").Value = "Source\N3XeS.CSharp.Common (Release)"
$ws.Range("C6").Value = "N3XeS.CSharp.Common.Exceptions"
$ws.Range("D6").Value = "ReadOnlyException"
$ws.Range("E6").Value = "ReadOnlyException(SerializationInfo, StreamingContext)"
$ws.Range("F6").Value = 95
$ws.Range("G6").Value = 1
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = 1

# Row 7 - ReadOnlyException(string)
$ws.Range("A7").Value = "Member"
$ws.Range("B7").Value = "Source\N3XeS.CSharp.Common (Release)"
$ws.Range("C7").Value = "N3XeS.CSharp.Common.Exceptions"
$ws.Range("D7").Value = "ReadOnlyException"
$ws.Range("E7").Value = "ReadOnlyException(string)"
$ws.Range("F7").Value = 95
$ws.Range("G7").Value = 1
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 1

# Row 8 - ReadOnlyException(string, Exception)
$ws.Range("A8").Value = "Member"
$ws.Range("B8").Value = "Source\N3XeS.CSharp.Common (Release)"
$ws.Range("C8").Value = "N3XeS.CSharp.Common.Exceptions"
$ws.Range("D8").Value = "ReadOnlyException"
$ws.Range("E8").Value = "ReadOnlyException(string, Exception)"
$ws.Range("F8").Value = 95
$ws.Range("G8").Value = 1
$ws.Range("I8").Value = 3
$ws.Range("J8").Value = 1

# Row 9 - N3XeS.CSharp.Common.Extensions (namespace)
$ws.Range("A9").Value = "Namespace"
$ws.Range("B9").Value = "Source\N3XeS.CSharp.Common (Release)"
$ws.Range("C9").Value = "N3XeS.CSharp.Common.Extensions"
$ws.Range("D9").Value = " "
$ws.Range("E9").Value = " "
$ws.Range("F9").Value = 92
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 4
$ws.Range("J9").Value = 1

# Row 10 - TypeExtension (type)
$ws.Range("A10").Value = "Type"
$ws.Range("B10").Value = "Source\N3XeS.CSharp.Common (Release)"
$ws.Range("C10").Value = "N3XeS.CSharp.Common.Extensions"
$ws.Range("D10").Value = "TypeExtension"
$ws.Range("E10").Value = " "
$ws.Range("F10").Value = 92
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 4
$ws.Range("J10").Value = 1

# Row 11 - GetActualType<T>(this T) : Type
$ws.Range("A11").Value = "Member"
$ws.Range("B11").Value = "Source\N3XeS.CSharp.Common (Release)"
$ws.Range("C11").Value = "N3XeS.CSharp.Common.Extensions"
$ws.Range("D11").Value = "TypeExtension"
$ws.Range("E11").Value = "GetActualType<T>(this T) : Type"
$ws.Range("F11").Value = 92
$ws.Range("G11").Value = 1
$ws.Range("I11").Value = 4
$ws.Range("J11").Value = 1

# Row 12 - N3XeS.CSharp.Common.Utilities (namespace)
$ws.Range("A12").Value = "Namespace"
$ws.Range("B12").Value = "Source\N3XeS.CSharp.Common (Release)"
$ws.Range("C12").Value = "N3XeS.CSharp.Common.Utilities"
$ws.Range("D12").Value = " "
$ws.Range("E12").Value = " "
$ws.Range("F12").Value = 74
$ws.Range("G12").Value = 2
$ws.Range("H12").Value = 1
$ws.Range("I12").Value = 3
$ws.Range("J12").Value = 4

# Row 13 - TypeUtility (type)
$ws.Range("A13").Value = "Type"
$ws.Range("B13").Value = "Source\N3XeS.CSharp.Common (Release)"
$ws.Range("C13").Value = "N3XeS.CSharp.Common.Utilities"
$ws.Range("D13").Value = "TypeUtility"
$ws.Range("E13").Value = " "
$ws.Range("F13").Value = 74
$ws.Range("G13").Value = 2
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 3
$ws.Range("J13").Value = 4

# Row 14 - GetActualType<T>(T) : Type
$ws.Range("A14").Value = "Member"
$ws.Range("B14").Value = "Source\N3XeS.CSharp.Common (Release)"
$ws.Range("C14").Value = "N3XeS.CSharp.Common.Utilities"
$ws.Range("D14").Value = "TypeUtility"
$ws.Range("E14").Value = "GetActualType<T>(T) : Type"
$ws.Range("F14").Value = 74
$ws.Range("G14").Value = 2
$ws.Range("I14").Value = 3
$ws.Range("J14").Value = 4

# ---------------------------------------------------------------------------
# Re-apply the AutoFilter so it (and the hidden _FilterDatabase defined
# name) spans the new A1:J14 extent instead of the old A1:J8.
# ---------------------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:J14").AutoFilter()

for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$J`$14"
    }
}

# ---------------------------------------------------------------------------
# Column widths were re-autofit by Excel after the data refresh (columns
# B-E shrink/grow to match the new longest entries).
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 38.7109375
$ws.Columns.Item(3).ColumnWidth = 32.85546875
$ws.Columns.Item(4).ColumnWidth = 18.5703125
$ws.Columns.Item(5).ColumnWidth = 52.7109375

# ---------------------------------------------------------------------------
# Mirror the final selection state recorded in the sheet (a "select all"
# click occurred after the refresh).
# ---------------------------------------------------------------------------
$ws.Cells.Select()
